# Applies the "Fix updating files in strike app" edits to card_1.docx:
#   - "2.5 Количество участников забастовки/акции" value:
#       51-100 человек -> менее 10 человек
#   - "3.3 Дата начало проведения забастовки/акции" value:
#       2021-08-19 18:00:00+00:00 -> 2021-08-28 18:00:00+00:00
#   - "3.4 Дата конца проведения забастовки/акции" value:
#       2021-08-19 18:00:00+00:00 -> 2021-08-24 18:00:00+00:00
#   - "3.8 Дата последних изменений" value:
#       2021-08-20 09:26:23.145947+00:00 -> 2021-08-25 10:58:07.159991+00:00
# (the sibling "3.7 Дата создания" row keeps its original timestamp and is
#  left untouched, even though its value text is identical before the edit)
#
# NOTE: the data-value cells repeat identical text across rows (e.g. both the
# start- and end-date rows read "2021-08-19 18:00:00+00:00" beforehand), so a
# document-wide Find/Replace (Replace:=wdReplaceAll) would touch every
# occurrence instead of just the intended one. Table.Cell(...).Range is also
# not reliable here for scoping a Find. Instead we locate the single
# paragraph that holds each value by its exact text and replace the text of
# just that paragraph's Range, built from precise Start/End offsets, which
# reliably scopes the edit to only that paragraph.

$d = $word.ActiveDocument

function Set-ParagraphText($oldText, $newText, $skip) {
    $seen = 0
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $p = $d.Paragraphs.Item($i)
        # Paragraph.Range.Text includes the trailing paragraph-mark (CR)
        # character, so trim it before comparing against the plain run text.
        $ptext = $p.Range.Text.TrimEnd([char]13)
        if ($ptext -eq $oldText) {
            $seen = $seen + 1
            if ($seen -gt $skip) {
                $r = $d.Range($p.Range.Start, $p.Range.End)
                # Replace:=1 (wdReplaceOne), NOT 2 (wdReplaceAll) -- with
                # ReplaceAll this runtime replaces every matching occurrence
                # in the whole document regardless of the Range used to
                # launch Find, even though the Range is scoped to a single
                # paragraph. wdReplaceOne correctly stays within the Range.
                $r.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                 $true, 1, $false, $newText, 1) | Out-Null
                return
            }
        }
    }
    throw "Paragraph with text '$oldText' (occurrence $($skip + 1)) not found"
}

# "2.5" row value (only occurrence of this text in the document).
Set-ParagraphText "51-100 человек" "менее 10 человек" 0

# "3.3" row value: first of the two identical "2021-08-19 ..." paragraphs.
Set-ParagraphText "2021-08-19 18:00:00+00:00" "2021-08-28 18:00:00+00:00" 0

# "3.4" row value: second of the two identical "2021-08-19 ..." paragraphs.
Set-ParagraphText "2021-08-19 18:00:00+00:00" "2021-08-24 18:00:00+00:00" 0

# "3.8" row value: second of the two identical "2021-08-20 09:26:23..."
# paragraphs (the first, "3.7" Дата создания, is left unchanged).
Set-ParagraphText "2021-08-20 09:26:23.145947+00:00" "2021-08-25 10:58:07.159991+00:00" 1
